$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells. Number/percent-looking cells are
# formatted as Text first so Excel keeps the original string formatting
# (e.g. "49.00", "0.0797", thousand-dot separators, padded percent signs)
# instead of silently converting them into floating point numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.071.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.302.37'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.87%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.43'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.95%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.30'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.04'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +16.46%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.662.15'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.372.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.807'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.962.57'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0907'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.11'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.91'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '237.49'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.20'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +14.61%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.47'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.59'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.80'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.89'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.18'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.54'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.07'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.78%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.67%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.82'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.37'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.998.48'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.10%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.14'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.78'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.97'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.527.92'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.60%  '
